$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "products"

$ws.Range("H1").Value = "images"
$ws.Range("I1").Value = "formatFix"
$ws.Range("F1").Value = "manage_stock"
$ws.Range("E1").Value = "slug"
$ws.Range("E2").Value = "test-product"
$ws.Range("G1").Value = "sku"
$ws.Range("G2").Value = "clitestproduct"

$ws.Range("D2").Value = 199
$ws.Range("F2").Value = $false

$ws.Columns.Item(6).ColumnWidth = 12.1
$ws.Columns.Item(7).ColumnWidth = 12.1

$ws.Range("G11").Select()
